# The commit adds a new (blank/zeroed) data row right under the header
# row of the "Đơn sale phụ" sheet, which expands the used range from
# A1:T1 to A1:T2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn sale phụ")

# Text columns on the new row are left blank.
$ws.Range("A2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("Q2").Value = ""
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = ""
$ws.Range("T2").Value = ""

# Numeric columns on the new row default to 0.
$ws.Range("B2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
